$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2789
$ws1.Range("F6").Value = 2440
$ws1.Range("F8").Value = 31
$ws1.Range("F9").Value = 33
$ws1.Range("F10").Value = 2879
$ws1.Range("F12").Value = 26
$ws1.Range("F13").Value = 6958
$ws1.Range("F14").Value = 270
$ws1.Range("F16").Value = 209
$ws1.Range("F18").Value = 468
$ws1.Range("F19").Value = 8064
$ws1.Range("F20").Value = 15
$ws1.Range("F22").Value = 256
$ws1.Range("F23").Value = 57
$ws1.Range("F27").Value = 64
$ws1.Range("F31").Value = 53
$ws1.Range("F33").Value = 2589
$ws1.Range("F38").Value = 86
$ws1.Range("F39").Value = 634
$ws1.Range("F40").Value = 3644
$ws1.Range("F41").Value = 163
$ws1.Range("F43").Value = 142

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 34
$ws2.Range("F4").Value = 26
$ws2.Range("F5").Value = 240

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 34
$ws4.Range("F3").Value = 2789
$ws4.Range("F4").Value = 26
$ws4.Range("F5").Value = 240
$ws4.Range("F6").Value = 240
$ws4.Range("F8").Value = 2440
$ws4.Range("F11").Value = 31
$ws4.Range("F12").Value = 33
$ws4.Range("F13").Value = 2879
$ws4.Range("F17").Value = 26
$ws4.Range("F18").Value = 6958
$ws4.Range("F19").Value = 270
$ws4.Range("F21").Value = 209
$ws4.Range("F23").Value = 468
$ws4.Range("F24").Value = 8064
$ws4.Range("F25").Value = 15
$ws4.Range("F27").Value = 256
$ws4.Range("F28").Value = 57
$ws4.Range("F32").Value = 64
$ws4.Range("F35").Value = 53
$ws4.Range("F38").Value = 2589
$ws4.Range("F43").Value = 86
$ws4.Range("F44").Value = 634
$ws4.Range("F46").Value = 3644
$ws4.Range("F47").Value = 163
$ws4.Range("F50").Value = 142
